# epexspot_prices.xlsx automated update
# - "Prix Spot" sheet: append a new day column (BF) with header "10-aug"
#   and 24 hourly prices.
# - "Gaz" sheet: append a new row (55) with date 2025-08-08 and its price.
# - "CO2" sheet: append a new row (55) with date 2025-08-08 and its price.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": new column BF (10-aug)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Header cell BF1 needs to look like the other header cells (bold,
# bordered, centered) -- copy the format from the previous header cell
# (BE1) before writing the new header text so we reuse the existing
# style instead of minting a new one.
$ws1.Range("BE1").Copy()
$ws1.Range("BF1").PasteSpecial(-4122)
$ws1.Cells.Item(1, 58).Value = "10-aug"

$ws1.Cells.Item(2, 58).Value = 77.88
$ws1.Cells.Item(3, 58).Value = 67.28
$ws1.Cells.Item(4, 58).Value = 58.05
$ws1.Cells.Item(5, 58).Value = 49.32
$ws1.Cells.Item(6, 58).Value = 46.22
$ws1.Cells.Item(7, 58).Value = 46.01
$ws1.Cells.Item(8, 58).Value = 42.41
$ws1.Cells.Item(9, 58).Value = 49.65
$ws1.Cells.Item(10, 58).Value = 26.87
$ws1.Cells.Item(11, 58).Value = -0.01
$ws1.Cells.Item(12, 58).Value = -2.26
$ws1.Cells.Item(13, 58).Value = -15.6
$ws1.Cells.Item(14, 58).Value = -14.05
$ws1.Cells.Item(15, 58).Value = -36.19
$ws1.Cells.Item(16, 58).Value = -50.29
$ws1.Cells.Item(17, 58).Value = -21
$ws1.Cells.Item(18, 58).Value = -1.16
$ws1.Cells.Item(19, 58).Value = 1.72
$ws1.Cells.Item(20, 58).Value = 26.5
$ws1.Cells.Item(21, 58).Value = 80.5
$ws1.Cells.Item(22, 58).Value = 98.01
$ws1.Cells.Item(23, 58).Value = 96.14
$ws1.Cells.Item(24, 58).Value = 90.59
$ws1.Cells.Item(25, 58).Value = 75.57

# ---------------------------------------------------------------------
# Sheet "Gaz": new row 55 (2025-08-08)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date column to stay plain text (matching the existing rows,
# which store ISO date strings as text, not date serials) by switching
# the cell to a text number format before assigning it, then restoring
# the default "Normal" style so no extra formatting is left behind.
$ws2.Cells.Item(55, 1).NumberFormat = "@"
$ws2.Cells.Item(55, 1).Value = "2025-08-08"
$ws2.Cells.Item(55, 1).Style = "Normal"
$ws2.Cells.Item(55, 2).Value = 31.2

# ---------------------------------------------------------------------
# Sheet "CO2": new row 55 (2025-08-08)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Cells.Item(55, 1).NumberFormat = "@"
$ws3.Cells.Item(55, 1).Value = "2025-08-08"
$ws3.Cells.Item(55, 1).Style = "Normal"
$ws3.Cells.Item(55, 2).Value = 71.75
